$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 9: A9 stays blank (but takes on the data row's number format/style), B9:W9 = 0
$ws.Range("A9:W9").NumberFormat = $ws.Range("A8:W8").NumberFormat

$ws.Range("B9:W9").Value = 0
